# Insert a new data row at row 5 (pushing the existing rows 5 and 6 down to
# rows 6 and 7), then populate the new row with the latest weekly price entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Insert()

$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44650
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 100112052
$ws.Range("G5").Value = "Albahaca"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 130
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = 3308
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 551
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = "Hortaliza"
